$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.446.07"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.107.01"
$ws.Range("E3").Value = "  +4.82%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'329.23"
$ws.Range("D6").Value = "'0.9986"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "'0.5265"
$ws.Range("E7").Value = "  +2.37%  "
$ws.Range("D8").Value = "'0.4353"
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("D9").Value = "'0.08849"
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("D10").Value = "'47.01"
$ws.Range("E10").Value = "  +8.71%  "
$ws.Range("D11").Value = "'1.165"
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("D12").Value = "'24.65"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "2.096.72"
$ws.Range("E13").Value = "  +4.29%  "
$ws.Range("D14").Value = "'6.728"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").Value = "'7.786"
$ws.Range("E15").Value = "  +4.48%  "
$ws.Range("D16").Value = "'96.37"
$ws.Range("E16").Value = "  +2.08%  "
$ws.Range("D17").Value = "'0.9984"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "'0.00001128"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").Value = "'0.06634"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "'0.9982"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").Value = "'6.345"
$ws.Range("E22").Value = "  +2.36%  "
$ws.Range("D23").Value = "30.507.14"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "'12.38"
$ws.Range("E24").Value = "  +5.05%  "
$ws.Range("D25").Value = "'2.326"
$ws.Range("E25").Value = "  +3.81%  "
$ws.Range("D26").Value = "2.338.79"
$ws.Range("E26").Value = "  +4.13%  "
$ws.Range("D27").Value = "'22.45"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").Value = "'2.592"
$ws.Range("E28").Value = "  +7.67%  "
$ws.Range("D29").Value = "'161.62"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "'132.53"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").Value = "'1.210"
$ws.Range("E31").Value = "  +6.71%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1075"
$ws.Range("E32").Value = "  +2.05%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'1.688"
$ws.Range("E33").Value = "  +24.00%  "
$ws.Range("D34").Value = "'6.184"
$ws.Range("E34").Value = "  +2.11%  "
$ws.Range("D35").Value = "'3.920"
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("D36").Value = "'9.959"
$ws.Range("E36").Value = "  +10.02%  "
$ws.Range("D37").Value = "'0.02582"
$ws.Range("E37").Value = "  +2.40%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.499"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.06711"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").Value = "'12.72"
$ws.Range("E40").Value = "  +3.09%  "
$ws.Range("E41").Value = "  +3.66%  "
$ws.Range("D42").Value = "'0.6834"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("D43").Value = "'1.257"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").Value = "'14.11"
$ws.Range("E44").Value = "  +3.84%  "
$ws.Range("D45").Value = "'0.9975"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("D46").Value = "'0.6382"
$ws.Range("E46").Value = "  +3.67%  "
$ws.Range("D47").Value = "'2.213"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").Value = "'3.616"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").Value = "'1.252"
$ws.Range("E49").Value = "  -0.78%  "
$ws.Range("D50").Value = "'1.199"
$ws.Range("E50").Value = "  +8.38%  "
$ws.Range("D51").Value = "'82.38"
$ws.Range("E51").Value = "  +2.27%  "
